$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A67").Value = "2024-10-06 00:00:00"
$ws.Range("B67").Value = 75650
$ws.Range("C67").Value = 10756.89
$ws.Range("D67").Value = 9519.370000000001
$ws.Range("E67").Value = 7.0184
